$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 31, pushing existing rows 31-45 down to 32-46.
$ws.Rows.Item(31).Insert()

# New weekly data point for the inserted row 31.
$ws.Cells.Item(31, 1).Value = 6
$ws.Cells.Item(31, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(31, 3).Value = "Metropolitana"
$ws.Cells.Item(31, 4).Value = 44806
$ws.Cells.Item(31, 4).NumberFormat = $ws.Cells.Item(32, 4).NumberFormat
$ws.Cells.Item(31, 5).Value = 13
$ws.Cells.Item(31, 6).Value = 100112035
$ws.Cells.Item(31, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 400
$ws.Cells.Item(31, 11).Value = 17000
$ws.Cells.Item(31, 12).Value = 18000
$ws.Cells.Item(31, 13).Value = 17425
$ws.Cells.Item(31, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(31, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(31, 16).Value = 1162
$ws.Cells.Item(31, 17).Value = 15
$ws.Cells.Item(31, 18).Value = "Hortaliza"
